$wb = $excel.ActiveWorkbook

# --- Add new sheet "2025-12-03" at the end of the workbook ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-12-03"

# Use the previous week's sheet as a formatting template (same header / highlight styles)
$wsPrev = $wb.Worksheets.Item("2025-11-26")

# --- Header row ---
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "volume"
$ws.Range("D1").Value = "publisher"

$wsPrev.Range("A1:D1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# --- Data rows: rank, title, volume ---
$data = @(
  @(1, 'ブルーロック', 36),
  @(2, 'GIANT KILLING', 68),
  @(3, '薬屋のひとりごと', 16),
  @(4, 'ブルーピリオド', 18),
  @(5, 'ドラフトキング', 23),
  @(6, 'ふたりソロキャンプ', 22),
  @(7, 'ベー革', 8),
  @(8, 'ゾン100~ゾンビになるまでにしたい100のこと~', 21),
  @(9, '超人X', 14),
  @(10, 'ワンパンマン', 35),
  @(11, 'スライムテイマーの異世界ライフ1', 1),
  @(12, 'ハズレ枠ので最強になった俺がすべてを蹂躙するまで', 13),
  @(13, '黄泉のツガイ', 11),
  @(14, 'MIX', 24),
  @(15, 'ジャンケットバンク', 20),
  @(16, 'とんがり帽子のアトリエ', 15),
  @(17, '金色のガッシュ!! 2 Page', 35),
  @(18, '捨てられ貴族の無人島のびのび開拓記~ようやく自由を手に入れたので、もふもふたちと気まぐれスローライフを満喫します~', 1),
  @(19, '山奥育ちの俺のゆるり異世界生活~もふもふと最強たちに可愛がられて、二度目の人生満喫中~', 1),
  @(20, 'カテナチオ', 8),
  @(21, '終末のワルキューレ禁伝 神々の黙示録', 4),
  @(22, 'オシバナ!', 1),
  @(23, 'ちいかわ なんか小さくてかわいいやつ', 8),
  @(24, '異世界でスローライフを(願望)', 10),
  @(25, '虚構推理', 24),
  @(26, '彼女、お借りします', 43),
  @(27, 'ガチアクタ', 17),
  @(28, '声が聞きたい小平くん', 1),
  @(29, 'おっさん異世界で最強になる', 1),
  @(30, '勇者パーティから追い出された不遇職、ユニークスキルで最強になる1', 1),
  @(31, 'うちの妻… 抱いてみませんか?', 1),
  @(32, '解雇された暗黒兵士のスローなセカンドライフ', 17),
  @(33, '魔女大戦 32人の異才の魔女は殺し合う', 14),
  @(34, 'みいちゃんと山田さん', 4),
  @(35, '出会って5秒でバトル', 29),
  @(36, 'KING GOLF', 43),
  @(37, 'DEAR BOYS ACT4', 21),
  @(38, '春の嵐とモンスター', 9),
  @(39, 'この度、公爵家の令嬢の婚約者となりました。しかし、噂では性格が悪く、十歳も年上です。 コミック版', 1),
  @(40, '復讐代行アプリ ~1タップで処刑します~', 7),
  @(41, '錬金術師の辺境再生スローライフ~S級パーティーで孤立した少女をかばったら追放されたので、一緒に幸せに暮らします~1', 1),
  @(42, 'ハズレ枠ので最強になった俺がすべてを蹂躙するまで', 1),
  @(43, 'みいちゃんと山田さん', 2),
  @(44, '蒼天の拳 リジェネシス', 6),
  @(45, '二階堂地獄ゴルフ', 10),
  @(46, 'みいちゃんと山田さん', 3),
  @(47, '龍とカメレオン', 8),
  @(48, 'ゆるキャン△', 18),
  @(49, '顔だけじゃ好きになりません', 16),
  @(50, '復讐代行アプリ ~1タップで処刑します~', 5),
  @(51, '復讐代行アプリ ~1タップで処刑します~', 6),
  @(52, 'おっさん異世界で最強になる', 2),
  @(53, 'おっさん異世界で最強になる', 3),
  @(54, '久々に健康診断を受けたら最強ステータスになっていた ~追放されたオッサン冒険者、今更英雄を目指す~', 1),
  @(55, '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが', 1),
  @(56, '魔法歌姫マジカルギンガ 第23話', 23),
  @(57, 'ギルド最弱のFランク冒険者、ユニークスキル「隠しクエスト」で最強に成り上がる1', 1),
  @(58, 'Sランク冒険者である俺の娘たちは重度のファザコンでした', 1),
  @(59, '異世界でスローライフを(願望)', 1),
  @(60, '終末のワルキューレ奇譚 ジャック・ザ・リッパーの事件簿', 9),
  @(61, '貸した魔力はで強制徴収~用済みとパーティー追放された俺は、可愛いサポート妖精と一緒に取り立てた魔力を運用して最強を目指す。~', 5),
  @(62, 'SPY×FAMILY', 16),
  @(63, '離婚予定の契約婚なのに、冷酷公爵様に執着されています', 3),
  @(64, 'ふたり鷹', 2),
  @(65, '今夜コインランドリーで逢いましょう', 1),
  @(66, '住所不定キャンパーはダンジョンでのんびりと暮らしたい', 1),
  @(67, '久々に健康診断を受けたら最強ステータスになっていた ~追放されたオッサン冒険者、今更英雄を目指す~', 2),
  @(68, '久々に健康診断を受けたら最強ステータスになっていた ~追放されたオッサン冒険者、今更英雄を目指す~', 3),
  @(69, '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが', 2),
  @(70, '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが', 3),
  @(71, '煙たい話', 1),
  @(72, 'イジワル同居人は御曹司!?', 1),
  @(73, '錬金術師の辺境再生スローライフ~S級パーティーで孤立した少女をかばったら追放されたので、一緒に幸せに暮らします~6', 6),
  @(74, '最強ギルドを追放された《植物王》、実は世界樹に選ばれていたので植物の力で無双します3', 3),
  @(75, '現実主義勇者の王国再建記I', 1),
  @(76, 'BLUE GIANT MOMENTUM', 6),
  @(77, '転生貴族の万能開拓~スキルを使っていたら最強領地になりました~', 13),
  @(78, '前田慶次 かぶき旅', 20),
  @(79, '国を追われた竜師さん、拾われた隣国でうっかり無双してしまう。~弱小国家が大陸最強の竜の楽園になるまで~', 4),
  @(80, '異世界らくらくサバイバル~生存スキル強者の俺が美少女四人と暮らす無人島生活~', 7),
  @(81, 'ボーイッシュ彼女が可愛すぎる 通常版', 4),
  @(82, '九条の大罪', 15),
  @(83, '星屑の王子様', 9),
  @(84, '釣って食べたいギャル澤さん', 4),
  @(85, '杖と剣のウィストリア', 14),
  @(86, 'ピンクとハバネロ', 14),
  @(87, 'ワンモア・テーブル!~人生やり直しレストラン~1', 1),
  @(88, '懲らしめて差し上げますっ!~おてんば王女の下剋上日記~', 1),
  @(89, '勇者は村人Aに恋をする', 1),
  @(90, 'これから幸せになります! 虐げられ令嬢ですが敵対国の公爵様に何故か溺愛されてます', 4),
  @(91, '今さらだけど異世界満喫! ~気づけばアラサー冒険者ですが、ゲーム知識で強くてニューゲーム~ 第1話-1', 1),
  @(92, 'ファントム無頼', 1),
  @(93, 'クリスティ・ハイテンション', 1),
  @(94, 'マウントするやつ黙らせます', 3),
  @(95, '今夜コインランドリーで逢いましょう', 2),
  @(96, '今夜コインランドリーで逢いましょう', 3),
  @(97, '住所不定キャンパーはダンジョンでのんびりと暮らしたい', 2),
  @(98, '住所不定キャンパーはダンジョンでのんびりと暮らしたい', 3),
  @(99, '世界にさよならのキスをして', 1),
  @(100, 'マウントするやつ黙らせます', 2)
)

$arr = New-Object 'object[,]' $data.Count,3
for ($i = 0; $i -lt $data.Count; $i++) {
    $arr[$i,0] = $data[$i][0]
    $arr[$i,1] = $data[$i][1]
    $arr[$i,2] = $data[$i][2]
}

$ws.Range("A2").Resize($data.Count, 3).Value = $arr

# --- Highlight style (fill) applied to specific "volume" cells (new-entry rows) ---
$styleRows = @(12, 19, 20, 23, 29, 30, 31, 32, 40, 42, 43, 44, 47, 53, 54, 55, 56, 58, 59, 60, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 75, 76, 88, 89, 90, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101)

$styleSrc = $wsPrev.Range("C21")
$styleSrc.Copy()
foreach ($r in $styleRows) {
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
}

$ws.Range("A1").Select()

Write-Output "Added sheet 2025-12-03 with $($data.Count) data rows and $($styleRows.Count) highlighted cells"
